$wb = $excel.ActiveWorkbook

# "Generate Report for Archive": the handoff status moved from
# "Ready for handoff" to "In Translation" on the Overview roll-up sheet
# (zh-cn/de-de status columns) and on each per-locale detail sheet. The
# shorter status text lets the status column narrow accordingly.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E:F").ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C:C").ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C:C").ColumnWidth = 12.5
